# "Updated symbol list" refresh of the cryptos sheet (coinranking.com scrape):
# the Price column (D) gets new quotes and a handful of rows' Coin/Link/
# Price/Volume(1h) cells (B:E) shift to reflect the updated coin ranking.
# Row 8 (MXToken) and the no-price "--" rows are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numeric-looking values ("245.68", "0.03234", ...)
# as literal text, not numbers. Pre-format D2:D50 as Text so the assignments
# below don't get auto-coerced to the Number type by Excel.
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = '245.68'
$ws.Range("D3").Value = '25.47'
$ws.Range("D4").Value = '5.133'
$ws.Range("D5").Value = '0.05590'
$ws.Range("D6").Value = '6.490'
$ws.Range("D7").Value = '3.028'
$ws.Range("D9").Value = '0.8512'
$ws.Range("D10").Value = '0.1340'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.06948'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = '0.03234'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.02858'
$ws.Range("E13").Value = '12BitrueCoinBTR'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09397'
$ws.Range("E14").Value = '13BitMartTokenBMX'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001523'
$ws.Range("E15").Value = '14BitForexTokenBF'

$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").Value = '0.0005957'
$ws.Range("E16").Value = '15OneONEWorstin24h'

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.006224'
$ws.Range("E17").Value = '16TigerCashTCH'

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.532'
$ws.Range("E18").Value = '17LEOLEO'

$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '2.118'
$ws.Range("E19").Value = '18BTSETokenBTSE'

$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '0.3168'
$ws.Range("E20").Value = '19BitpandaEcosystemTokenBEST'

$ws.Range("D21").Value = '0.1319'
$ws.Range("D22").Value = '3.751'
$ws.Range("D23").Value = '0.04687'
$ws.Range("D24").Value = '0.1374'
$ws.Range("D25").Value = '0.001251'
$ws.Range("D26").Value = '0.004605'

$ws.Range("D27").Value = '0.00009596'
$ws.Range("E27").Value = '26NitroExNTXBestin24h'

$ws.Range("D28").Value = '0.0001389'
$ws.Range("D40").Value = '0.03658'

$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1054'
$ws.Range("E41").Value = '40BKEXTokenBKK'

$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '0.006117'
$ws.Range("E42").Value = '41KickTokenKICK'

$ws.Range("D43").Value = '0.002476'
$ws.Range("D44").Value = '0.007390'
$ws.Range("D45").Value = '0.00005319'

$ws.Range("D47").Value = '0.1334'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'

$ws.Range("D49").Value = '0.00002099'
$ws.Range("D50").Value = '0.0001999'
